# Correção ao use case "Confecionar Receita"
#
# The underlying text of a handful of use-case steps was rewritten, and the
# selection / scroll position of the sheet was moved further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (1/2).2 / Alternativa 1 - "ação do passo" -> "passo" ; still about the selected step
$ws.Cells.Item(14, 2).Value = " Alternativa 1 [passo suscita dúvidas] (Passo 1 e 2)"
$ws.Cells.Item(15, 4).Value = "(1/2).2 <<include>> Apresenta informação extra relativa ao passo selecionado"

# Alternativa 2 -> Exceção 2 (aspect of final result doesn't match the picture)
$ws.Cells.Item(18, 2).Value = " Exceção 2 [aspecto final não corresponde à imagem apresentada] (Passo 4)"

# 4.3 now returns to "Preparar Receita" instead of looping back to step 1
$ws.Cells.Item(20, 4).Value = "4.3 Regressa a Preparar Receita"

# Alternativa 3 now references Passo 5 (not Passo 4), and its response step is renumbered 5.1
$ws.Cells.Item(22, 2).Value = " Alternativa 3 [existem mais passos] (Passo 5)"
$ws.Cells.Item(22, 4).Value = "5.1 Regressa a 1"

# Move the viewport / selection further down the sheet (was topLeftCell A7 / D21)
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D23").Select() | Out-Null
